$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recursive update to sbol2 assembly: re-order which part slot each
# SynBioHub part URL is assigned to (swap A5<->B2 and B1<->B3 slots,
# i.e. cells E3/H3 and G3/I3 swap their displayed URL text).
$ws.Range("E3").Value = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/B0015/1"
$ws.Range("G3").Value = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/J23100/1"
$ws.Range("H3").Value = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/B0032/1"
$ws.Range("I3").Value = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/E0040m_gfp/1"
